$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.053.82"
$ws.Range("E2").Value = "  -2.38%  "
$ws.Range("D3").Value = "1.823.11"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("D4").Value = "'" + "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.10%  "
$ws.Range("D5").Value = "'" + "311.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.48%  "
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("D7").Value = "'" + "0.4225"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.77%  "
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("D9").Value = "'" + "0.07221"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.58%  "
$ws.Range("D10").Value = "'" + "0.8415"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.00%  "
$ws.Range("E11").Value = "  -3.69%  "
$ws.Range("D12").Value = "1.822.46"
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("D14").Value = "'" + "0.07067"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").Value = "'" + "5.294"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.57%  "
$ws.Range("D16").Value = "'" + "90.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("D17").Value = "'" + "1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("D18").Value = "'" + "0.000008744"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.64%  "
$ws.Range("D19").Value = "'" + "1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("D20").Value = "'" + "14.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.11%  "
$ws.Range("D21").Value = "27.162.00"
$ws.Range("E21").Value = "  -2.02%  "
$ws.Range("D22").Value = "'" + "5.141"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.16%  "
$ws.Range("D23").Value = "'" + "10.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.87%  "
$ws.Range("D24").Value = "2.054.16"
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("D25").Value = "'" + "1.982"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "'" + "151.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.12%  "
$ws.Range("D27").Value = "'" + "2.252"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.01%  "
$ws.Range("D28").Value = "'" + "18.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.05%  "
$ws.Range("D29").Value = "'" + "5.255"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.06%  "
$ws.Range("D30").Value = "'" + "117.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.09%  "
$ws.Range("D31").Value = "'" + "0.08716"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("D32").Value = "'" + "1.179"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.74%  "
$ws.Range("E33").Value = "  -5.09%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'" + "4.422"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.73%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'" + "2.883"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("D36").Value = "'" + "1.001"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("D37").Value = "'" + "1.088"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.60%  "
$ws.Range("D38").Value = "'" + "0.01950"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.65%  "
$ws.Range("D39").Value = "'" + "0.05254"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("D40").Value = "'" + "7.325"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("D41").Value = "'" + "2.869"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.71%  "
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D43").Value = "'" + "0.5042"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.67%  "
$ws.Range("D44").Value = "'" + "8.549"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.04%  "
$ws.Range("D45").Value = "'" + "10.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("D46").Value = "'" + "106.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.48%  "
$ws.Range("D47").Value = "'" + "0.4709"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.32%  "
$ws.Range("D48").Value = "'" + "1.921"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.51%  "
$ws.Range("D50").Value = "'" + "0.06335"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("D51").Value = "'" + "1.650"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.18%  "
